# Merge the split "<id>...</id>" runs back into a single run for each
# of the three newly-downloaded tc/tcn/tl entries (p100v_2, p101r_1,
# p101r_2). The "fig_p101r_1" entry is intentionally left untouched.

$d = $word.ActiveDocument

$ids = @("p100v_2", "p101r_1", "p101r_2")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $needle, 2) | Out-Null
}
